# Sample Project / Main.xlsx - "SAVE" edit
#
# The rules table on the active sheet has row 11 (the "R40" rule row).
# Cell B11 previously held the text "R40"; it is changed to hold the
# text "1" (still a text value, not a number).
#
# A leading apostrophe forces Excel to store the new value as text
# (t="s", i.e. a shared string) instead of re-interpreting "1" as a
# number, which mirrors exactly how this edit would be made by typing
# into the cell in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
